# Refresh cached Universalis market-price figures in the Durandal_Profits workbook.
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
#  LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns, per scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 4835925
$ws.Range("I28").Value = 6541933.5
$ws.Range("K28").Value = 6541933.5
$ws.Range("M28").Value = -6541448.5
# Row 33
$ws.Range("H33").Value = 6151.222
$ws.Range("I33").Value = 66.666664
$ws.Range("J33").Value = 7368.1333
$ws.Range("K33").Value = 66.666664
$ws.Range("L33").Value = 7368.1333
$ws.Range("M33").Value = 162.333336
$ws.Range("N33").Value = -7826.1333
# Row 41
$ws.Range("H41").Value = 6173510.5
$ws.Range("I41").Value = 15873288
$ws.Range("J41").Value = 925.36365
$ws.Range("K41").Value = 15873288
$ws.Range("L41").Value = 925.36365
$ws.Range("M41").Value = -15872848
$ws.Range("N41").Value = -1805.36365
# Row 62
$ws.Range("H62").Value = 2520.55
$ws.Range("I62").Value = 2431.5625
$ws.Range("J62").Value = 2876.5
$ws.Range("K62").Value = 2431.5625
$ws.Range("L62").Value = 2876.5
$ws.Range("M62").Value = -1807.5625
$ws.Range("N62").Value = -4124.5
# Row 65
$ws.Range("H65").Value = 2520.55
$ws.Range("I65").Value = 2431.5625
$ws.Range("J65").Value = 2876.5
$ws.Range("K65").Value = 12157.8125
$ws.Range("L65").Value = 14382.5
$ws.Range("M65").Value = -9037.8125
$ws.Range("N65").Value = -20622.5
# Row 76
$ws.Range("H76").Value = 2648747.5
$ws.Range("I76").Value = 3706646.2
$ws.Range("J76").Value = 4001
$ws.Range("K76").Value = 3706646.2
$ws.Range("L76").Value = 4001
$ws.Range("M76").Value = -3706331.2
$ws.Range("N76").Value = -4631
# Row 79
$ws.Range("H79").Value = 2648747.5
$ws.Range("I79").Value = 3706646.2
$ws.Range("J79").Value = 4001
$ws.Range("K79").Value = 3706646.2
$ws.Range("L79").Value = 4001
$ws.Range("M79").Value = -3705554.2
$ws.Range("N79").Value = -6185
# Row 82
$ws.Range("H82").Value = 6750
$ws.Range("I82").Value = 4733.3335
$ws.Range("K82").Value = 14200.0005
$ws.Range("M82").Value = -13794.0005
# Row 85
$ws.Range("H85").Value = 6750
$ws.Range("I85").Value = 4733.3335
$ws.Range("K85").Value = 14200.0005
$ws.Range("M85").Value = -12796.0005
# Row 86
$ws.Range("H86").Value = 1730.8966
$ws.Range("I86").Value = 1636
$ws.Range("J86").Value = 2029.1428
$ws.Range("K86").Value = 1636
$ws.Range("L86").Value = 2029.1428
$ws.Range("M86").Value = -513
$ws.Range("N86").Value = -4275.1428
# Row 89
$ws.Range("H89").Value = 1730.8966
$ws.Range("I89").Value = 1636
$ws.Range("J89").Value = 2029.1428
$ws.Range("K89").Value = 8180
$ws.Range("L89").Value = 10145.714
$ws.Range("M89").Value = -2564
$ws.Range("N89").Value = -21377.714
# Row 92
$ws.Range("H92").Value = 687.8823
$ws.Range("I92").Value = 734.1539
$ws.Range("J92").Value = 537.5
$ws.Range("K92").Value = 734.1539
$ws.Range("L92").Value = 537.5
$ws.Range("M92").Value = 513.8461
$ws.Range("N92").Value = -3033.5
# Row 98
$ws.Range("H98").Value = 65420372
$ws.Range("I98").Value = 9343026
$ws.Range("J98").Value = 200006000
$ws.Range("K98").Value = 9343026
$ws.Range("L98").Value = 200006000
$ws.Range("M98").Value = -9341528
$ws.Range("N98").Value = -200008996
# Row 106
$ws.Range("H106").Value = 2169.5217
$ws.Range("I106").Value = 2050.4736
$ws.Range("J106").Value = 2735
$ws.Range("K106").Value = 2050.4736
$ws.Range("L106").Value = 2735
$ws.Range("M106").Value = -1419.4736
$ws.Range("N106").Value = -3997
# Row 107
$ws.Range("H107").Value = 259.6087
$ws.Range("I107").Value = 257.77777
$ws.Range("J107").Value = 266.2
$ws.Range("K107").Value = 257.77777
$ws.Range("L107").Value = 266.2
$ws.Range("M107").Value = 1662.22223
$ws.Range("N107").Value = -4106.2
# Row 122
$ws.Range("H122").Value = 65420372
$ws.Range("I122").Value = 9343026
$ws.Range("J122").Value = 200006000
$ws.Range("K122").Value = 28029078
$ws.Range("L122").Value = 600018000
$ws.Range("M122").Value = -28026628
$ws.Range("N122").Value = -600022900
# Row 134
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 479406.84
$ws.Range("I32").Value = 2754.3774
$ws.Range("K32").Value = 2754.3774
$ws.Range("M32").Value = -2467.3774
# Row 41
$ws.Range("H41").Value = 1500
$ws.Range("I41").Value = 1500
$ws.Range("K41").Value = 1500
$ws.Range("M41").Value = -1086
# Row 131
$ws.Range("H131").Value = 64987.8
$ws.Range("J131").Value = 64987.8
$ws.Range("L131").Value = 64987.8
$ws.Range("N131").Value = -75067.8

$ws = $wb.Worksheets.Item("CRP")
# Row 141
$ws.Range("H141").Value = 22957.916
$ws.Range("I141").Value = 17500
$ws.Range("J141").Value = 23454.092
$ws.Range("K141").Value = 17500
$ws.Range("L141").Value = 23454.092
$ws.Range("M141").Value = -12320
$ws.Range("N141").Value = -33814.092

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1256.5
$ws.Range("I5").Value = 1258.6666
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 3775.9998
$ws.Range("L5").Value = 3750
$ws.Range("M5").Value = -3663.9998
$ws.Range("N5").Value = -3974
# Row 131
$ws.Range("H131").Value = 31251450
$ws.Range("I131").Value = 937.8
$ws.Range("J131").Value = 83335630
$ws.Range("K131").Value = 2813.4
$ws.Range("L131").Value = 250006890
$ws.Range("M131").Value = 2226.6
$ws.Range("N131").Value = -250016970
# Row 135
$ws.Range("H135").Value = 1256.5
$ws.Range("I135").Value = 1258.6666
$ws.Range("J135").Value = 1250
$ws.Range("K135").Value = 11327.9994
$ws.Range("L135").Value = 11250
$ws.Range("M135").Value = -8792.999400000001
$ws.Range("N135").Value = -16320

$ws = $wb.Worksheets.Item("GSM")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
# Row 128
$ws.Range("H128").Value = 45992
$ws.Range("J128").Value = 45992
$ws.Range("L128").Value = 45992
$ws.Range("N128").Value = -55952

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
# Row 68
$ws.Range("H68").Value = 1787.4773
$ws.Range("I68").Value = 1764.5186
$ws.Range("J68").Value = 1823.9412
$ws.Range("K68").Value = 1764.5186
$ws.Range("L68").Value = 1823.9412
$ws.Range("M68").Value = -1015.5186
$ws.Range("N68").Value = -3321.9412
# Row 71
$ws.Range("H71").Value = 1787.4773
$ws.Range("I71").Value = 1764.5186
$ws.Range("J71").Value = 1823.9412
$ws.Range("K71").Value = 8822.593000000001
$ws.Range("L71").Value = 9119.706
$ws.Range("M71").Value = -5078.593000000001
$ws.Range("N71").Value = -16607.706
